$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (17:20 -> 17:50)
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 17:50"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 174684
$ws.Range("C4").Value = 10896
$ws.Range("D4").Value = 6210
$ws.Range("E4").Value = 165074
$ws.Range("F4").Value = 3893
$ws.Range("G4").Value = 259
$ws.Range("H4").Value = 3400

# Row 21: Brasil
$ws.Range("A21").Value = "Brasil"
$ws.Range("B21").Value = 4685
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 127
$ws.Range("E21").Value = 4390
$ws.Range("F21").Value = 296
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 168

# Row 22: Noruega
$ws.Range("A22").Value = "Noruega"
$ws.Range("B22").Value = 4605
$ws.Range("C22").Value = 160
$ws.Range("D22").Value = 13
$ws.Range("E22").Value = 4556
$ws.Range("F22").Value = 97
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 36

# Row 27: Dinamarca
$ws.Range("A27").Value = "Dinamarca"
$ws.Range("B27").Value = 2860
$ws.Range("C27").Value = 283
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 2769
$ws.Range("F27").Value = 145
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 90

# Row 42: Sudafrica
$ws.Range("A42").Value = "Sudafrica"
$ws.Range("B42").Value = 1353
$ws.Range("C42").Value = 27
$ws.Range("D42").Value = 31
$ws.Range("E42").Value = 1319
$ws.Range("F42").Value = 7
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 3

# Row 43: Grecia
$ws.Range("A43").Value = "Grecia"
$ws.Range("B43").Value = 1314
$ws.Range("C43").Value = 102
$ws.Range("D43").Value = 52
$ws.Range("E43").Value = 1213
$ws.Range("F43").Value = 72
$ws.Range("G43").Value = 3
$ws.Range("H43").Value = 49

# Row 44: India
$ws.Range("A44").Value = "India"
$ws.Range("B44").Value = 1251
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 102
$ws.Range("E44").Value = 1117
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 32

# Row 49: Argentina
$ws.Range("A49").Value = "Argentina"
$ws.Range("B49").Value = 966
$ws.Range("C49").Value = 146
$ws.Range("D49").Value = 240
$ws.Range("E49").Value = 700
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 26

# Row 79: Kazajistan
$ws.Range("A79").Value = "Kazajistan"
$ws.Range("B79").Value = 340
$ws.Range("C79").Value = 38
$ws.Range("D79").Value = 22
$ws.Range("E79").Value = 316
$ws.Range("F79").Value = 6
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 2

# Row 81: Republica de Macedonia
$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("B81").Value = 329
$ws.Range("C81").Value = 44
$ws.Range("D81").Value = 12
$ws.Range("E81").Value = 308
$ws.Range("F81").Value = 4
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 9

# Row 139: Togo
$ws.Range("A139").Value = "Togo"
$ws.Range("B139").Value = 34
$ws.Range("C139").Value = 4
$ws.Range("D139").Value = 10
$ws.Range("E139").Value = 23
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 1

# Row 140: Uganda
$ws.Range("A140").Value = "Uganda"
$ws.Range("B140").Value = 33
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 0
$ws.Range("E140").Value = 33
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 0

# Row 141: El Salvador
$ws.Range("A141").Value = "El Salvador"
$ws.Range("B141").Value = 32
$ws.Range("C141").Value = 2
$ws.Range("D141").Value = 0
$ws.Range("E141").Value = 32
$ws.Range("F141").Value = 5
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 0

# Row 142: Guam
$ws.Range("A142").Value = "Guam"
$ws.Range("B142").Value = 32
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 0
$ws.Range("E142").Value = 31
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 1

# Row 143: Republica de Yibuti
$ws.Range("A143").Value = "Republica de Yibuti"
$ws.Range("B143").Value = 30
$ws.Range("C143").Value = 12
$ws.Range("D143").Value = 0
$ws.Range("E143").Value = 30
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0

# Row 157: Birmania
$ws.Range("A157").Value = "Birmania"
$ws.Range("B157").Value = 15
$ws.Range("C157").Value = 1
$ws.Range("D157").Value = 0
$ws.Range("E157").Value = 14
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 1

# Row 158: San Martin (Parte Francesa)
$ws.Range("A158").Value = "San Martin (Parte Francesa)"
$ws.Range("B158").Value = 15
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 2
$ws.Range("E158").Value = 12
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 1

# Row 159: Bahamas
$ws.Range("A159").Value = "Bahamas"
$ws.Range("B159").Value = 14
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 1
$ws.Range("E159").Value = 13
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 0
